# Generate Report for Handback
#
# Semantics (derived from the target diff):
#  - Overview sheet: the "Ready for handoff" status for both language rows
#    becomes "Handed back: in sync with en-US" (shared by E2/F2/E3/F3 since
#    they all shared the same status string).
#  - zh-cn / de-de sheets: the "Latest Target File" (col I) and
#    "Latest Handback File" (col J) columns get populated for both rows, I
#    getting a hyperlink (same target url as the row's source-file link in
#    col A), and the "Latest Handback DateTime" (col K) gets a real
#    timestamp instead of the 0001-01-01 placeholder (zh-cn finished a
#    little before de-de).
#  - A handful of columns get widened so the new long filename values are
#    readable.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: Excel's ColumnWidth setter rounds through an integer pixel grid
# (pixels = floor(chars*6 + 5 + 0.5); stored width = pixels/6), so to land
# on a particular *stored* width we have to solve for the "chars" input
# that snaps to the nearest matching pixel bucket.
# ---------------------------------------------------------------------
function Set-StoredColumnWidth($col, [double]$targetStoredWidth) {
    $mdw = 6.0
    $pad = 5.0
    $px = [Math]::Round($targetStoredWidth * $mdw)
    $lo = ($px - $pad - 0.5) / $mdw
    $hi = ($px + 1 - $pad - 0.5) / $mdw
    $col.ColumnWidth = ($lo + $hi) / 2.0
}

# ---------------------------------------------------------------------
# 1) Overview sheet: handback status text + widened status columns
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

Set-StoredColumnWidth $overview.Columns.Item(5) 29.9777047293527
Set-StoredColumnWidth $overview.Columns.Item(6) 29.9777047293527

# ---------------------------------------------------------------------
# 2) Per-language sheets (zh-cn, de-de): fill in Latest Target File /
#    Latest Handback File / Latest Handback DateTime for both rows.
# ---------------------------------------------------------------------
function Update-LanguageSheet($sheetName, $handbackFile2, $handbackFile3, $handbackDateTime) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Column widths: Status (C), Latest Target File (I), Latest Handback File (J)
    Set-StoredColumnWidth $ws.Columns.Item(3) 29.9777047293527
    Set-StoredColumnWidth $ws.Columns.Item(9) 40
    Set-StoredColumnWidth $ws.Columns.Item(10) 40

    # Look up the existing hyperlinks on column A (source file) so the new
    # "Latest Target File" links point at the same targets.
    $rowLinks = @{}
    foreach ($h in $ws.Hyperlinks) {
        $rowLinks[$h.Range.Row] = @{ Address = $h.Address; Display = $h.TextToDisplay }
    }

    # Row 2
    $link2 = $rowLinks[2]
    $ws.Hyperlinks.Add($ws.Range("I2"), $link2.Address, "", "", $link2.Display) | Out-Null
    $ws.Range("J2").Value = $handbackFile2
    $ws.Range("K2").Value = $handbackDateTime

    # Row 3
    $link3 = $rowLinks[3]
    $ws.Hyperlinks.Add($ws.Range("I3"), $link3.Address, "", "", $link3.Display) | Out-Null
    $ws.Range("J3").Value = $handbackFile3
    $ws.Range("K3").Value = $handbackDateTime
}

Update-LanguageSheet "zh-cn" `
    "4e151b7d-c896-4d62-befc-fa85d0c7fb64.e91b9791dcc307957f144a3d2fa2ac8419832c38.zh-cn.xlf" `
    "c7c04534-ef70-484b-9a08-6b57324de0e7.734b7f93500b0a2d5f6b43fda6b08856f7709149.zh-cn.xlf" `
    "2016-08-23 06:45:37"

Update-LanguageSheet "de-de" `
    "4e151b7d-c896-4d62-befc-fa85d0c7fb64.e91b9791dcc307957f144a3d2fa2ac8419832c38.de-de.xlf" `
    "c7c04534-ef70-484b-9a08-6b57324de0e7.734b7f93500b0a2d5f6b43fda6b08856f7709149.de-de.xlf" `
    "2016-08-23 06:45:44"

Write-Output "Handback report generated"
